$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.62005599999999
$ws.Range("H2").Value = 289.860168
$ws.Range("I2").Value = 0.2116037895476247
$ws.Range("J2").Value = 0.2183905833651517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.757653333333334
$ws.Range("N2").Value = 11.27296
$ws.Range("O2").Value = 0.0006253057195333491
$ws.Range("P2").Value = 0.0006253596992176285
$ws.Range("Q2").Value = 363.0646754952533
$ws.Range("R2").Value = 3267.58207945728
$ws.Range("S2").Value = 0.0001323170598790608
$ws.Range("T2").Value = 0.0001365726695251937
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.62005599999999
$ws.Range("H3").Value = 289.860168
$ws.Range("I3").Value = 0.2116037895476247
$ws.Range("J3").Value = 0.2183905833651517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.337854
$ws.Range("N3").Value = 19.013562
$ws.Range("O3").Value = 0.001054673224006999
$ws.Range("P3").Value = 0.001054764268956488
$ws.Range("Q3").Value = 612.3638083998239
$ws.Range("R3").Value = 5511.274275598415
$ws.Range("S3").Value = 0.0002231728509342919
$ws.Range("T3").Value = 0.0002303505840101252
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.62005599999999
$ws.Range("H4").Value = 289.860168
$ws.Range("I4").Value = 0.2116037895476247
$ws.Range("J4").Value = 0.2183905833651517
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3733.712707666667
$ws.Range("N4").Value = 11201.138123
$ws.Range("O4").Value = 0.6213217942399283
$ws.Range("P4").Value = 0.6213754300107861
$ws.Range("Q4").Value = 360751.530902665
$ws.Range("R4").Value = 3246763.778123985
$ws.Range("S4").Value = 0.1314740461896983
$ws.Range("T4").Value = 0.1357025426488276
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.62005599999999
$ws.Range("H5").Value = 289.860168
$ws.Range("I5").Value = 0.2116037895476247
$ws.Range("J5").Value = 0.2183905833651517
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2263.941528333333
$ws.Range("N5").Value = 6791.824585
$ws.Range("O5").Value = 0.3767392733645568
$ws.Range("P5").Value = 0.3767717954835726
$ws.Range("Q5").Value = 218742.1572482922
$ws.Range("R5").Value = 1968679.41523463
$ws.Range("S5").Value = 0.07971945791535871
$ws.Range("T5").Value = 0.08228341221119305
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 96.62005599999999
$ws.Range("H6").Value = 289.860168
$ws.Range("I6").Value = 0.2116037895476247
$ws.Range("J6").Value = 0.2183905833651517
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 1.5561305
$ws.Range("N6").Value = 3.112261
$ws.Range("O6").Value = 0.0002589534519745364
$ws.Range("P6").Value = 0.0001726505374672451
$ws.Range("Q6").Value = 150.353416053308
$ws.Range("R6").Value = 902.120496319848
$ws.Range("S6").Value = 0.00005479553175425074
$ws.Range("T6").Value = 0.00003770525159577865
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 122.3539896666667
$ws.Range("H7").Value = 367.061969
$ws.Range("I7").Value = 0.2679626668787852
$ws.Range("J7").Value = 0.2765570657541026
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.757653333333334
$ws.Range("N7").Value = 11.27296
$ws.Range("O7").Value = 0.0006253057195333491
$ws.Range("P7").Value = 0.0006253596992176285
$ws.Range("Q7").Value = 459.7638771175822
$ws.Range("R7").Value = 4137.87489405824
$ws.Range("S7").Value = 0.0001675585882207139
$ws.Range("T7").Value = 0.0001729476434564955
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 122.3539896666667
$ws.Range("H8").Value = 367.061969
$ws.Range("I8").Value = 0.2679626668787852
$ws.Range("J8").Value = 0.2765570657541026
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.337854
$ws.Range("N8").Value = 19.013562
$ws.Range("O8").Value = 0.001054673224006999
$ws.Range("P8").Value = 0.001054764268956488
$ws.Range("Q8").Value = 775.461722824842
$ws.Range("R8").Value = 6979.155505423578
$ws.Range("S8").Value = 0.000282613049790562
$ws.Range("T8").Value = 0.0002917025112848774
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 122.3539896666667
$ws.Range("H9").Value = 367.061969
$ws.Range("I9").Value = 0.2679626668787852
$ws.Range("J9").Value = 0.2765570657541026
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3733.712707666667
$ws.Range("N9").Value = 11201.138123
$ws.Range("O9").Value = 0.6213217942399283
$ws.Range("P9").Value = 0.6213754300107861
$ws.Range("Q9").Value = 456834.6460521494
$ws.Range("R9").Value = 4111511.814469344
$ws.Range("S9").Value = 0.166491044974443
$ws.Range("T9").Value = 0.1718457656554767
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.3539896666667
$ws.Range("H10").Value = 367.061969
$ws.Range("I10").Value = 0.2679626668787852
$ws.Range("J10").Value = 0.2765570657541026
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2263.941528333333
$ws.Range("N10").Value = 6791.824585
$ws.Range("O10").Value = 0.3767392733645568
$ws.Range("P10").Value = 0.3767717954835726
$ws.Range("Q10").Value = 277002.2783636342
$ws.Range("R10").Value = 2493020.505272708
$ws.Range("S10").Value = 0.1009520604087423
$ws.Range("T10").Value = 0.1041989022178417
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.3539896666667
$ws.Range("H11").Value = 367.061969
$ws.Range("I11").Value = 0.2679626668787852
$ws.Range("J11").Value = 0.2765570657541026
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 1.5561305
$ws.Range("N11").Value = 3.112261
$ws.Range("O11").Value = 0.0002589534519745364
$ws.Range("P11").Value = 0.0001726505374672451
$ws.Range("Q11").Value = 190.3987751169848
$ws.Range("R11").Value = 1142.392650701909
$ws.Range("S11").Value = 0.00006938985758856422
$ws.Range("T11").Value = 0.00004774772604281007
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 83.74384566666667
$ws.Range("H12").Value = 251.231537
$ws.Range("I12").Value = 0.1834041070557659
$ws.Range("J12").Value = 0.1892864490617203
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.757653333333334
$ws.Range("N12").Value = 11.27296
$ws.Range("O12").Value = 0.0006253057195333491
$ws.Range("P12").Value = 0.0006253596992176285
$ws.Range("Q12").Value = 314.6803408155023
$ws.Range("R12").Value = 2832.12306733952
$ws.Range("S12").Value = 0.0001146836371278771
$ws.Range("T12").Value = 0.0001183721168512104
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 83.74384566666667
$ws.Range("H13").Value = 251.231537
$ws.Range("I13").Value = 0.1834041070557659
$ws.Range("J13").Value = 0.1892864490617203
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.337854
$ws.Range("N13").Value = 19.013562
$ws.Range("O13").Value = 0.001054673224006999
$ws.Range("P13").Value = 0.001054764268956488
$ws.Range("Q13").Value = 530.7562672338661
$ws.Range("R13").Value = 4776.806405104794
$ws.Range("S13").Value = 0.0001934314008846294
$ws.Range("T13").Value = 0.0001996525830679549
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 83.74384566666667
$ws.Range("H14").Value = 251.231537
$ws.Range("I14").Value = 0.1834041070557659
$ws.Range("J14").Value = 0.1892864490617203
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3733.712707666667
$ws.Range("N14").Value = 11201.138123
$ws.Range("O14").Value = 0.6213217942399283
$ws.Range("P14").Value = 0.6213754300107861
$ws.Range("Q14").Value = 312675.4607545095
$ws.Range("R14").Value = 2814079.146790585
$ws.Range("S14").Value = 0.1139529688668603
$ws.Range("T14").Value = 0.1176179486809412
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 83.74384566666667
$ws.Range("H15").Value = 251.231537
$ws.Range("I15").Value = 0.1834041070557659
$ws.Range("J15").Value = 0.1892864490617203
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2263.941528333333
$ws.Range("N15").Value = 6791.824585
$ws.Range("O15").Value = 0.3767392733645568
$ws.Range("P15").Value = 0.3767717954835726
$ws.Range("Q15").Value = 189591.1699471042
$ws.Range("R15").Value = 1706320.529523937
$ws.Range("S15").Value = 0.06909553002426462
$ws.Range("T15").Value = 0.07131779527369418
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 83.74384566666667
$ws.Range("H16").Value = 251.231537
$ws.Range("I16").Value = 0.1834041070557659
$ws.Range("J16").Value = 0.1892864490617203
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 1.5561305
$ws.Range("N16").Value = 3.112261
$ws.Range("O16").Value = 0.0002589534519745364
$ws.Range("P16").Value = 0.0001726505374672451
$ws.Range("Q16").Value = 130.3163524291928
$ws.Range("R16").Value = 781.898114575157
$ws.Range("S16").Value = 0.00004749312662839801
$ws.Range("T16").Value = 0.00003268040716577234
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 111.321218
$ws.Range("H17").Value = 333.963654
$ws.Range("I17").Value = 0.2438002270031519
$ws.Range("J17").Value = 0.2516196610353779
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.757653333333334
$ws.Range("N17").Value = 11.27296
$ws.Range("O17").Value = 0.0006253057195333491
$ws.Range("P17").Value = 0.0006253596992176285
$ws.Range("Q17").Value = 418.3065458884267
$ws.Range("R17").Value = 3764.758912995841
$ws.Range("S17").Value = 0.0001524496763685997
$ws.Range("T17").Value = 0.0001573527955423255
# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 111.321218
$ws.Range("H18").Value = 333.963654
$ws.Range("I18").Value = 0.2438002270031519
$ws.Range("J18").Value = 0.2516196610353779
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 6.337854
$ws.Range("N18").Value = 19.013562
$ws.Range("O18").Value = 0.001054673224006999
$ws.Range("P18").Value = 0.001054764268956488
$ws.Range("Q18").Value = 705.537626786172
$ws.Range("R18").Value = 6349.838641075548
$ws.Range("S18").Value = 0.0002571295714270525
$ws.Range("T18").Value = 0.0002653994278270596
# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 111.321218
$ws.Range("H19").Value = 333.963654
$ws.Range("I19").Value = 0.2438002270031519
$ws.Range("J19").Value = 0.2516196610353779
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 3733.712707666667
$ws.Range("N19").Value = 11201.138123
$ws.Range("O19").Value = 0.6213217942399283
$ws.Range("P19").Value = 0.6213754300107861
$ws.Range("Q19").Value = 415641.4462795313
$ws.Range("R19").Value = 3740773.016515782
$ws.Range("S19").Value = 0.1514783944777001
$ws.Range("T19").Value = 0.1563502750750262
# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 111.321218
$ws.Range("H20").Value = 333.963654
$ws.Range("I20").Value = 0.2438002270031519
$ws.Range("J20").Value = 0.2516196610353779
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 2263.941528333333
$ws.Range("N20").Value = 6791.824585
$ws.Range("O20").Value = 0.3767392733645568
$ws.Range("P20").Value = 0.3767717954835726
$ws.Range("Q20").Value = 252024.7284148482
$ws.Range("R20").Value = 2268222.555733634
$ws.Range("S20").Value = 0.09184912036728143
$ws.Range("T20").Value = 0.09480319146726725
# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 111.321218
$ws.Range("H21").Value = 333.963654
$ws.Range("I21").Value = 0.2438002270031519
$ws.Range("J21").Value = 0.2516196610353779
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 1.5561305
$ws.Range("N21").Value = 3.112261
$ws.Range("O21").Value = 0.0002589534519745364
$ws.Range("P21").Value = 0.0001726505374672451
$ws.Range("Q21").Value = 173.230342626949
$ws.Range("R21").Value = 1039.382055761694
$ws.Range("S21").Value = 0.00006313291037464177
$ws.Range("T21").Value = 0.00004344226971508402
# Row 22
$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 42.569235
$ws.Range("H22").Value = 85.13847
$ws.Range("I22").Value = 0.09322920951467238
$ws.Range("J22").Value = 0.06414624078364733
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 3.757653333333334
$ws.Range("N22").Value = 11.27296
$ws.Range("O22").Value = 0.0006253057195333491
$ws.Range("P22").Value = 0.0006253596992176285
$ws.Range("Q22").Value = 159.9604277952
$ws.Range("R22").Value = 959.7625667712001
$ws.Range("S22").Value = 0.00005829675793709756
$ws.Range("T22").Value = 0.00004011447384240327
# Row 23
$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 42.569235
$ws.Range("H23").Value = 85.13847
$ws.Range("I23").Value = 0.09322920951467238
$ws.Range("J23").Value = 0.06414624078364733
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 6.337854
$ws.Range("N23").Value = 19.013562
$ws.Range("O23").Value = 0.001054673224006999
$ws.Range("P23").Value = 0.001054764268956488
$ws.Range("Q23").Value = 269.79759632169
$ws.Range("R23").Value = 1618.78557793014
$ws.Range("S23").Value = 0.00009832635097046353
$ws.Range("T23").Value = 0.00006765916276647063
# Row 24
$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 42.569235
$ws.Range("H24").Value = 85.13847
$ws.Range("I24").Value = 0.09322920951467238
$ws.Range("J24").Value = 0.06414624078364733
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 3733.712707666667
$ws.Range("N24").Value = 11201.138123
$ws.Range("O24").Value = 0.6213217942399283
$ws.Range("P24").Value = 0.6213754300107861
$ws.Range("Q24").Value = 158941.2936751487
$ws.Range("R24").Value = 953647.7620508919
$ws.Range("S24").Value = 0.05792533973122644
$ws.Range("T24").Value = 0.03985889795051429
# Row 25
$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 42.569235
$ws.Range("H25").Value = 85.13847
$ws.Range("I25").Value = 0.09322920951467238
$ws.Range("J25").Value = 0.06414624078364733
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 2263.941528333333
$ws.Range("N25").Value = 6791.824585
$ws.Range("O25").Value = 0.3767392733645568
$ws.Range("P25").Value = 0.3767717954835726
$ws.Range("Q25").Value = 96374.25894588082
$ws.Range("R25").Value = 578245.553675285
$ws.Range("S25").Value = 0.0351231046489097
$ws.Range("T25").Value = 0.02416849431357638
# Row 26
$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 42.569235
$ws.Range("H26").Value = 85.13847
$ws.Range("I26").Value = 0.09322920951467238
$ws.Range("J26").Value = 0.06414624078364733
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 1.5561305
$ws.Range("N26").Value = 3.112261
$ws.Range("O26").Value = 0.0002589534519745364
$ws.Range("P26").Value = 0.0001726505374672451
$ws.Range("Q26").Value = 66.2432849451675
$ws.Range("R26").Value = 264.97313978067
$ws.Range("S26").Value = 0.00002414202562868171
$ws.Range("T26").Value = 0.00001107488294780003
